$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("attributes")

# Insert a new column at H ("description"), shifting H..Q to I..R.
$ws.Columns.Item(8).Insert()

# Header for the new column.
$ws.Cells.Item(1, 8).Value2 = "description"

# Fill the description column with a copy of the attribute name (column A)
# for every data row, mirroring the fixture's existing name/description
# pairing used on the "entities" sheet.
for ($r = 2; $r -le 52; $r++) {
    $name = $ws.Cells.Item($r, 1).Value2
    if ($name -ne $null) {
        $ws.Cells.Item($r, 8).Value2 = $name
    }
}

# Match the column width Excel computed for the new text column.
$ws.Columns.Item(8).ColumnWidth = 28.14

# Restore the active selection Excel left on the new column.
[void]$ws.Range("H5").Select()
